$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 50,4
$data[0,0] = 0.1068081706762314
$data[0,1] = 0.9673090577125549
$data[0,2] = 0.0221733134239912
$data[0,3] = 0.9959924221038818
$data[1,0] = 0.02469750680029392
$data[1,1] = 0.9943469762802124
$data[1,2] = 0.00997348316013813
$data[1,3] = 0.9988341331481934
$data[2,0] = 0.0159869696944952
$data[2,1] = 0.995252251625061
$data[2,2] = 0.006918943952769041
$data[2,3] = 0.9975954294204712
$data[3,0] = 0.012208366766572
$data[3,1] = 0.9950108528137207
$data[3,2] = 0.004813738167285919
$data[3,3] = 0.998105525970459
$data[4,0] = 0.01113185100257397
$data[4,1] = 0.995252251625061
$data[4,2] = 0.001581466873176396
$data[4,3] = 0.9993442296981812
$data[5,0] = 0.01062070205807686
$data[5,1] = 0.9951717853546143
$data[5,2] = 0.001875111483968794
$data[5,3] = 0.9993442296981812
$data[6,0] = 0.01133355498313904
$data[6,1] = 0.9949706196784973
$data[6,2] = 0.001390560530126095
$data[6,3] = 0.999489963054657
$data[7,0] = 0.009661545045673847
$data[7,1] = 0.9956948757171631
$data[7,2] = 0.001344242598861456
$data[7,3] = 0.999489963054657
$data[8,0] = 0.009431255050003529
$data[8,1] = 0.9954534769058228
$data[8,2] = 0.001576627953909338
$data[8,3] = 0.9991984963417053
$data[9,0] = 0.009200803004205227
$data[9,1] = 0.9954735636711121
$data[9,2] = 0.001212380127981305
$data[9,3] = 0.999489963054657
$data[10,0] = 0.01015620399266481
$data[10,1] = 0.9950712323188782
$data[10,2] = 0.002772331004962325
$data[10,3] = 0.9991984963417053
$data[11,0] = 0.009123739786446095
$data[11,1] = 0.9956546425819397
$data[11,2] = 0.001736307633109391
$data[11,3] = 0.9993442296981812
$data[12,0] = 0.009336655959486961
$data[12,1] = 0.9951919317245483
$data[12,2] = 0.002222855109721422
$data[12,3] = 0.9993442296981812
$data[13,0] = 0.008759592659771442
$data[13,1] = 0.9958155751228333
$data[13,2] = 0.001724094850942492
$data[13,3] = 0.9991984963417053
$data[14,0] = 0.00934179313480854
$data[14,1] = 0.9954333305358887
$data[14,2] = 0.001524470513686538
$data[14,3] = 0.9993442296981812
$data[15,0] = 0.009614747017621994
$data[15,1] = 0.9953126311302185
$data[15,2] = 0.004663704894483089
$data[15,3] = 0.9993442296981812
$data[16,0] = 0.009268315508961678
$data[16,1] = 0.9951315522193909
$data[16,2] = 0.002417960204184055
$data[16,3] = 0.9993442296981812
$data[17,0] = 0.008870344609022141
$data[17,1] = 0.9954534769058228
$data[17,2] = 0.001171692390926182
$data[17,3] = 0.999489963054657
$data[18,0] = 0.009318721480667591
$data[18,1] = 0.9955741763114929
$data[18,2] = 0.001172867137938738
$data[18,3] = 0.999489963054657
$data[19,0] = 0.008593840524554253
$data[19,1] = 0.9958758950233459
$data[19,2] = 0.001350579201243818
$data[19,3] = 0.9993442296981812
$data[20,0] = 0.009307581000030041
$data[20,1] = 0.9955339431762695
$data[20,2] = 0.001290653832256794
$data[20,3] = 0.999489963054657
$data[21,0] = 0.009420525282621384
$data[21,1] = 0.9952120184898376
$data[21,2] = 0.001219181809574366
$data[21,3] = 0.999489963054657
$data[22,0] = 0.008851123042404652
$data[22,1] = 0.9955741763114929
$data[22,2] = 0.002010358031839132
$data[22,3] = 0.9993442296981812
$data[23,0] = 0.009741540066897869
$data[23,1] = 0.9948901534080505
$data[23,2] = 0.00380669254809618
$data[23,3] = 0.9993442296981812
$data[24,0] = 0.01077430881559849
$data[24,1] = 0.9949706196784973
$data[24,2] = 0.001256832503713667
$data[24,3] = 0.999489963054657
$data[25,0] = 0.009092814289033413
$data[25,1] = 0.995674729347229
$data[25,2] = 0.001407375908456743
$data[25,3] = 0.9994170665740967
$data[26,0] = 0.008677861653268337
$data[26,1] = 0.9955339431762695
$data[26,2] = 0.001287330058403313
$data[26,3] = 0.9993442296981812
$data[27,0] = 0.008799066767096519
$data[27,1] = 0.9956546425819397
$data[27,2] = 0.001282678917050362
$data[27,3] = 0.9993442296981812
$data[28,0] = 0.0091822799295187
$data[28,1] = 0.9951114654541016
$data[28,2] = 0.0009748386219143867
$data[28,3] = 0.9994170665740967
$data[29,0] = 0.008524461649358273
$data[29,1] = 0.9954333305358887
$data[29,2] = 0.001098708482459188
$data[29,3] = 0.9993442296981812
$data[30,0] = 0.00949336402118206
$data[30,1] = 0.9953126311302185
$data[30,2] = 0.00107584975194186
$data[30,3] = 0.9993442296981812
$data[31,0] = 0.009108642116189003
$data[31,1] = 0.9949907660484314
$data[31,2] = 0.00102001300547272
$data[31,3] = 0.999489963054657
$data[32,0] = 0.008581042289733887
$data[32,1] = 0.9957551956176758
$data[32,2] = 0.001042782212607563
$data[32,3] = 0.9994170665740967
$data[33,0] = 0.009838566184043884
$data[33,1] = 0.9950108528137207
$data[33,2] = 0.00976625457406044
$data[33,3] = 0.9992713332176208
$data[34,0] = 0.008432622998952866
$data[34,1] = 0.9957753419876099
$data[34,2] = 0.001126538380049169
$data[34,3] = 0.9993442296981812
$data[35,0] = 0.008704917505383492
$data[35,1] = 0.9953930974006653
$data[35,2] = 0.0009408654295839369
$data[35,3] = 0.999489963054657
$data[36,0] = 0.008937020786106586
$data[36,1] = 0.9954131841659546
$data[36,2] = 0.008867413736879826
$data[36,3] = 0.9992713332176208
$data[37,0] = 0.008721551857888699
$data[37,1] = 0.9950712323188782
$data[37,2] = 0.00340975821018219
$data[37,3] = 0.9993442296981812
$data[38,0] = 0.00842635240405798
$data[38,1] = 0.9954333305358887
$data[38,2] = 0.003168723313137889
$data[38,3] = 0.9993442296981812
$data[39,0] = 0.008780546486377716
$data[39,1] = 0.9955339431762695
$data[39,2] = 0.002084123436361551
$data[39,3] = 0.9991984963417053
$data[40,0] = 0.008719589561223984
$data[40,1] = 0.9951315522193909
$data[40,2] = 0.002670306479558349
$data[40,3] = 0.9991984963417053
$data[41,0] = 0.00840804073959589
$data[41,1] = 0.9954131841659546
$data[41,2] = 0.004714268259704113
$data[41,3] = 0.9991984963417053
$data[42,0] = 0.00871859397739172
$data[42,1] = 0.9956144094467163
$data[42,2] = 0.002476723166182637
$data[42,3] = 0.9991984963417053
$data[43,0] = 0.008898387663066387
$data[43,1] = 0.9954534769058228
$data[43,2] = 0.004443750716745853
$data[43,3] = 0.9990527629852295
$data[44,0] = 0.008342387154698372
$data[44,1] = 0.9958356618881226
$data[44,2] = 0.003080145921558142
$data[44,3] = 0.9990527629852295
$data[45,0] = 0.008390046656131744
$data[45,1] = 0.9958356618881226
$data[45,2] = 0.003876763628795743
$data[45,3] = 0.9991984963417053
$data[46,0] = 0.009328382089734077
$data[46,1] = 0.9948901534080505
$data[46,2] = 0.001664511626586318
$data[46,3] = 0.9991984963417053
$data[47,0] = 0.008551171980798244
$data[47,1] = 0.9956144094467163
$data[47,2] = 0.003393363440409303
$data[47,3] = 0.9991984963417053
$data[48,0] = 0.008087413385510445
$data[48,1] = 0.9957753419876099
$data[48,2] = 0.004807431250810623
$data[48,3] = 0.9978868961334229
$data[49,0] = 0.008895776234567165
$data[49,1] = 0.9953528642654419
$data[49,2] = 0.002282810397446156
$data[49,3] = 0.9993442296981812

$ws.Range("A2:D51").Value = $data
